$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 19.93587785526113
$ws.Range("C2").Value = 9.474960569113012
$ws.Range("D2").Value = 8.10468717712218
$ws.Range("E2").Value = 12.69313455404209
$ws.Range("F2").Value = 35.67475656164123
$ws.Range("J2").Value = 9.93129353983379
$ws.Range("L2").Value = 11.21381394625832
$ws.Range("O2").Value = 27.4364470073285
$ws.Range("B3").Value = 19.39582765267483
$ws.Range("C3").Value = 9.188252085129255
$ws.Range("D3").Value = 8.09812369119679
$ws.Range("E3").Value = 12.7300669480721
$ws.Range("F3").Value = 35.81446367997764
$ws.Range("J3").Value = 9.959825859776213
$ws.Range("L3").Value = 11.18387958347973
$ws.Range("O3").Value = 27.56367958688861
$ws.Range("B4").Value = 19.05820643016793
$ws.Range("C4").Value = 9.00695041679678
$ws.Range("D4").Value = 8.094827478781148
$ws.Range("E4").Value = 12.75427584390294
$ws.Range("F4").Value = 35.90987989271568
$ws.Range("J4").Value = 9.978246839644353
$ws.Range("L4").Value = 11.16662915140338
$ws.Range("O4").Value = 27.64887083926283
$ws.Range("B5").Value = 18.91929841816435
$ws.Range("C5").Value = 8.931830301088867
$ws.Range("D5").Value = 8.093669891693935
$ws.Range("E5").Value = 12.76452704375937
$ws.Range("F5").Value = 35.95117724214575
$ws.Range("J5").Value = 9.985981044286991
$ws.Range("L5").Value = 11.15988749693787
$ws.Range("O5").Value = 27.6853594949701
$ws.Range("B6").Value = 18.89615866105433
$ws.Range("C6").Value = 8.919284499285784
$ws.Range("D6").Value = 8.093488919649563
$ws.Range("E6").Value = 12.76625257090624
$ws.Range("F6").Value = 35.9581802114649
$ws.Range("J6").Value = 9.987279064026442
$ws.Range("L6").Value = 11.15878556121454
$ws.Range("O6").Value = 27.69152531028189
$ws.Range("B7").Value = 19.05633817026191
$ws.Range("C7").Value = 9.005942216709279
$ws.Range("D7").Value = 8.094811114026271
$ws.Range("E7").Value = 12.75441253198791
$ws.Range("F7").Value = 35.91042707852537
$ws.Range("J7").Value = 9.978350223762982
$ws.Range("L7").Value = 11.16653705965082
$ws.Range("O7").Value = 27.64935576815446
$ws.Range("B8").Value = 19.75103075546535
$ws.Range("C8").Value = 9.377244163334725
$ws.Range("D8").Value = 8.10227247933682
$ws.Range("E8").Value = 12.70555117263828
$ws.Range("F8").Value = 35.72092272928638
$ws.Range("J8").Value = 9.940944683314559
$ws.Range("L8").Value = 11.20326065523017
$ws.Range("O8").Value = 27.47884628519461
$ws.Range("B9").Value = 21.0577030607798
$ws.Range("C9").Value = 10.06025832788181
$ws.Range("D9").Value = 8.122675681328474
$ws.Range("E9").Value = 12.62186580809294
$ws.Range("F9").Value = 35.42612554349309
$ws.Range("J9").Value = 9.874718186583802
$ws.Range("L9").Value = 11.28405212049566
$ws.Range("O9").Value = 27.20080800569284
$ws.Range("B10").Value = 21.97427599352037
$ws.Range("C10").Value = 10.53056907621581
$ws.Range("D10").Value = 8.141114242718215
$ws.Range("E10").Value = 12.56774175984621
$ws.Range("F10").Value = 35.25685513862346
$ws.Range("J10").Value = 9.83036124945763
$ws.Range("L10").Value = 11.34850287659259
$ws.Range("O10").Value = 27.03119049047493
$ws.Range("B11").Value = 22.38005833456704
$ws.Range("C11").Value = 10.73697490172286
$ws.Range("D11").Value = 8.150235116588801
$ws.Range("E11").Value = 12.54470980996095
$ws.Range("F11").Value = 35.19022169840531
$ws.Range("J11").Value = 9.811106317825862
$ws.Range("L11").Value = 11.37887062419716
$ws.Range("O11").Value = 26.9616196010302
$ws.Range("B12").Value = 22.5319817165591
$ws.Range("C12").Value = 10.81400007320556
$ws.Range("D12").Value = 8.153792882682142
$ws.Range("E12").Value = 12.5362161969064
$ws.Range("F12").Value = 35.16648782045977
$ws.Range("J12").Value = 9.803947030009223
$ws.Range("L12").Value = 11.3905156392658
$ws.Range("O12").Value = 26.93637170920278
$ws.Range("B13").Value = 22.49934153505555
$ws.Range("C13").Value = 10.79746262307884
$ws.Range("D13").Value = 8.153022058547887
$ws.Range("E13").Value = 12.53803531149311
$ws.Range("F13").Value = 35.17153255836234
$ws.Range("J13").Value = 9.805483043843754
$ws.Range("L13").Value = 11.38800129761122
$ws.Range("O13").Value = 26.94176040920662
$ws.Range("B14").Value = 22.39259269388036
$ws.Range("C14").Value = 10.74333485725698
$ws.Range("D14").Value = 8.150525743010272
$ws.Range("E14").Value = 12.54400646672028
$ws.Range("F14").Value = 35.1882390174657
$ws.Range("J14").Value = 9.810514674262048
$ws.Range("L14").Value = 11.37982578392721
$ws.Range("O14").Value = 26.95952042137388
$ws.Range("B15").Value = 22.32697597961226
$ws.Range("C15").Value = 10.71003062635929
$ws.Range("D15").Value = 8.149010160764036
$ws.Range("E15").Value = 12.54769366177015
$ws.Range("F15").Value = 35.1986676129005
$ws.Range("J15").Value = 9.813613882779059
$ws.Range("L15").Value = 11.37483681901117
$ws.Range("O15").Value = 26.97054198677944
$ws.Range("B16").Value = 21.94752087255077
$ws.Range("C16").Value = 10.51692366078785
$ws.Range("D16").Value = 8.140532794442407
$ws.Range("E16").Value = 12.56927889423798
$ws.Range("F16").Value = 35.2614191057517
$ws.Range("J16").Value = 9.831638129198012
$ws.Range("L16").Value = 11.34653891864664
$ws.Range("O16").Value = 27.03589031951676
$ws.Range("B17").Value = 21.71177894044974
$ws.Range("C17").Value = 10.3964886133869
$ws.Range("D17").Value = 8.135518896727344
$ws.Range("E17").Value = 12.58292747640883
$ws.Range("F17").Value = 35.30257639969734
$ws.Range("J17").Value = 9.842931442350723
$ws.Range("L17").Value = 11.32944389687418
$ws.Range("O17").Value = 27.07792727926891
$ws.Range("B18").Value = 21.5751419595811
$ws.Range("C18").Value = 10.32651048857633
$ws.Range("D18").Value = 8.132704131509751
$ws.Range("E18").Value = 12.59092740419647
$ws.Range("F18").Value = 35.3272244792256
$ws.Range("J18").Value = 9.849513996285495
$ws.Range("L18").Value = 11.31971034230999
$ws.Range("O18").Value = 27.1028197932573
$ws.Range("B19").Value = 21.52870383960412
$ws.Range("C19").Value = 10.30269733443677
$ws.Range("D19").Value = 8.131763013234524
$ws.Range("E19").Value = 12.59366175455104
$ws.Range("F19").Value = 35.33573719381204
$ws.Range("J19").Value = 9.85175768697499
$ws.Range("L19").Value = 11.31643190186328
$ws.Range("O19").Value = 27.11137041704165
$ws.Range("B20").Value = 21.73698315834119
$ws.Range("C20").Value = 10.40938270318714
$ws.Range("D20").Value = 8.136045494387441
$ws.Range("E20").Value = 12.58145907933873
$ws.Range("F20").Value = 35.29809411801507
$ws.Range("J20").Value = 9.841720256588628
$ws.Range("L20").Value = 11.33125347847744
$ws.Range("O20").Value = 27.07337843791536
$ws.Range("B21").Value = 22.4239955578805
$ws.Range("C21").Value = 10.75926471585169
$ws.Range("D21").Value = 8.151256164338418
$ws.Range("E21").Value = 12.54224640711777
$ws.Range("F21").Value = 35.18329118952621
$ws.Range("J21").Value = 9.809033181254543
$ws.Range("L21").Value = 11.38222322800761
$ws.Range("O21").Value = 26.95427405455399
$ws.Range("B22").Value = 22.86281765702426
$ws.Range("C22").Value = 10.98128993499753
$ws.Range("D22").Value = 8.161802011738841
$ws.Range("E22").Value = 12.51794791324244
$ws.Range("F22").Value = 35.11700008018918
$ws.Range("J22").Value = 9.788440195462824
$ws.Range("L22").Value = 11.41637983770337
$ws.Range("O22").Value = 26.88282979775409
$ws.Range("B23").Value = 22.629581705507
$ws.Range("C23").Value = 10.86341437870539
$ws.Range("D23").Value = 8.156118675392115
$ws.Range("E23").Value = 12.53079500007771
$ws.Range("F23").Value = 35.1515787858581
$ws.Range("J23").Value = 9.799360818851545
$ws.Range("L23").Value = 11.39807432735969
$ws.Range("O23").Value = 26.92037368598078
$ws.Range("B24").Value = 21.72559176566823
$ws.Range("C24").Value = 10.40355558122436
$ws.Range("D24").Value = 8.135807208323662
$ws.Range("E24").Value = 12.58212246448974
$ws.Range("F24").Value = 35.30011748613903
$ws.Range("J24").Value = 9.842267553654322
$ws.Range("L24").Value = 11.33043507156161
$ws.Range("O24").Value = 27.07543271143611
$ws.Range("B25").Value = 20.711202308535
$ws.Range("C25").Value = 9.88075375985502
$ws.Range("D25").Value = 8.116545193993085
$ws.Range("E25").Value = 12.64320990978359
$ws.Range("F25").Value = 35.49760049959088
$ws.Range("J25").Value = 9.891875957499948
$ws.Range("L25").Value = 11.26128207412268
$ws.Range("O25").Value = 27.269961289148
